# Refresh the cryptos price table (columns D = Price, E = Volume(1h))
# with the latest scraped values. Numeric-looking "Price" strings are
# written through a temporary Text number format so Excel keeps them as
# literal strings (e.g. "25.31" stays text instead of becoming 25.31 as
# a number); ClearFormats() afterwards restores the original (default)
# cell style so only the displayed text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = '29.829.77'
$ws.Cells.Item(2, 5).Value = '  -0.25%  '
$ws.Cells.Item(3, 4).Value = '1.887.02'
$ws.Cells.Item(3, 5).Value = '  -0.42%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
Set-TextValue 5 4 '0.7517'
$ws.Cells.Item(5, 5).Value = '  -3.67%  '
$ws.Cells.Item(6, 5).Value = '  -0.65%  '
Set-TextValue 8 4 '0.3117'
$ws.Cells.Item(8, 5).Value = '  -0.43%  '
Set-TextValue 9 4 '25.31'
$ws.Cells.Item(9, 5).Value = '  -1.51%  '
Set-TextValue 10 4 '0.07116'
$ws.Cells.Item(10, 5).Value = '  -2.98%  '
Set-TextValue 11 4 '0.08473'
$ws.Cells.Item(11, 5).Value = '  +4.76%  '
Set-TextValue 12 4 '0.7595'
$ws.Cells.Item(12, 5).Value = '  -1.64%  '
$ws.Cells.Item(13, 4).Value = '1.881.33'
$ws.Cells.Item(13, 5).Value = '  -0.84%  '
Set-TextValue 14 4 '5.357'
$ws.Cells.Item(14, 5).Value = '  -2.62%  '
Set-TextValue 15 4 '93.32'
$ws.Cells.Item(15, 5).Value = '  -0.60%  '
Set-TextValue 16 4 '6.135'
$ws.Cells.Item(16, 5).Value = '  -1.43%  '
$ws.Cells.Item(17, 4).Value = '29.756.48'
$ws.Cells.Item(17, 5).Value = '  -0.27%  '
Set-TextValue 18 4 '13.70'
$ws.Cells.Item(18, 5).Value = '  -1.85%  '
Set-TextValue 19 4 '243.61'
$ws.Cells.Item(19, 5).Value = '  -1.54%  '
Set-TextValue 20 4 '0.000007805'
$ws.Cells.Item(20, 5).Value = '  -0.21%  '
Set-TextValue 21 4 '0.9996'
$ws.Cells.Item(21, 5).Value = '  +0.02%  '
$ws.Cells.Item(22, 4).Value = '2.137.33'
$ws.Cells.Item(22, 5).Value = '  +0.97%  '
Set-TextValue 23 4 '7.991'
$ws.Cells.Item(23, 5).Value = '  -1.35%  '
$ws.Cells.Item(24, 5).Value = '  +0.00%  '
Set-TextValue 25 4 '0.1593'
$ws.Cells.Item(25, 5).Value = '  +0.13%  '
Set-TextValue 26 4 '9.365'
$ws.Cells.Item(26, 5).Value = '  -0.86%  '
Set-TextValue 27 4 '162.88'
$ws.Cells.Item(27, 5).Value = '  -0.41%  '
Set-TextValue 28 4 '18.72'
$ws.Cells.Item(28, 5).Value = '  +0.00%  '
$ws.Cells.Item(29, 5).Value = '  +0.21%  '
Set-TextValue 30 4 '1.491'
$ws.Cells.Item(30, 5).Value = '  +3.85%  '
Set-TextValue 31 4 '1.537'
$ws.Cells.Item(31, 5).Value = '  -0.53%  '
Set-TextValue 32 4 '4.510'
$ws.Cells.Item(32, 5).Value = '  +0.54%  '
Set-TextValue 33 4 '4.125'
$ws.Cells.Item(33, 5).Value = '  +1.61%  '
Set-TextValue 34 4 '0.05418'
$ws.Cells.Item(35, 5).Value = '  +0.22%  '
Set-TextValue 36 4 '0.7499'
$ws.Cells.Item(36, 5).Value = '  -0.37%  '
Set-TextValue 37 4 '1.001'
$ws.Cells.Item(37, 5).Value = '  -0.37%  '
$ws.Cells.Item(38, 5).Value = '  +1.04%  '
Set-TextValue 39 4 '0.01946'
$ws.Cells.Item(39, 5).Value = '  +0.62%  '
Set-TextValue 40 4 '2.773'
$ws.Cells.Item(40, 5).Value = '  -0.87%  '
Set-TextValue 41 4 '0.4453'
$ws.Cells.Item(41, 5).Value = '  -0.39%  '
Set-TextValue 42 4 '6.099'
$ws.Cells.Item(42, 5).Value = '  +2.28%  '
$ws.Cells.Item(43, 4).Value = '1.091.12'
$ws.Cells.Item(43, 5).Value = '  -1.91%  '
Set-TextValue 44 4 '72.52'
$ws.Cells.Item(44, 5).Value = '  -2.36%  '
Set-TextValue 45 4 '0.8592'
$ws.Cells.Item(45, 5).Value = '  +0.91%  '
$ws.Cells.Item(46, 5).Value = '  +0.02%  '
Set-TextValue 47 4 '7.727'
$ws.Cells.Item(47, 5).Value = '  +2.77%  '
Set-TextValue 48 4 '102.27'
$ws.Cells.Item(48, 5).Value = '  -0.24%  '
Set-TextValue 49 4 '1.858'
$ws.Cells.Item(49, 5).Value = '  -1.53%  '
Set-TextValue 50 4 '3.060'
$ws.Cells.Item(50, 5).Value = '  +0.49%  '
$ws.Cells.Item(51, 4).Value = '2.038.37'
$ws.Cells.Item(51, 5).Value = '  -0.23%  '
